$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to make room for the "Version" column.
$ws.Columns("A:A").Insert()

# Header
$ws.Range("A1").Value = "Version"

# Data rows - constant version value "1.0" for every existing data row
$ws.Range("A2").Value = "1.0"
$ws.Range("A3").Value = "1.0"
$ws.Range("A4").Value = "1.0"
$ws.Range("A5").Value = "1.0"
$ws.Range("A6").Value = "1.0"
